$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(58).Insert()
$ws.Range("A57:Q57").Copy()
$ws.Range("A58:Q58").PasteSpecial(-4122)
$ws.Range("C58").Value2
